$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    # Use a narrow Find (no in-line replacement argument) so the engine
    # does not run its "replace straight quotes with smart quotes"
    # AutoFormat pass on the replacement text; a bare Find.Execute just
    # collapses/extends $r to the matched span, which we can then set
    # directly via .Text (this also removes any w:proofErr markers that
    # are fully inside the matched span, keeping spellStart/spellEnd and
    # gramStart/gramEnd balanced).
    $r = $d.Content
    $ok = $r.Find.Execute($find)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $find"
        return
    }
    $r.Text = $replace
}

# --- Paragraph 1 (About me intro) ---
Replace-Text "girlfriend of 2 years" "girlfriend of two years"
Replace-Text "also almost my neighbor, so I have the 2 people closest" "also almost my neighbour, so I have the two people closest"
Replace-Text "closest to me at arms length. I got" "closest to me at arm’s length. I got"
Replace-Text "had a programming Course, and now" "had a programming course, and now"
Replace-Text "pursuing and Education in Computer Science" "pursuing an Education in Computer Science"

# --- Paragraph 2 (club / PlayStation) ---
Replace-Text 'get into the "club" in my school' 'get into the "club" at my school'
Replace-Text "which all the kids in the school was fighting" "which all the kids in the school were fighting"
Replace-Text "we played a lot together, it even was used" "we played a lot together. It even was used"
Replace-Text "to often hold friendly(Very very serious and competitive) family tournaments, and even since then, my love for videogames was cemented." "to often hold friendly (and by friendly, I of course mean very very serious and competitive) family tournaments. This love for videogames became even stronger when my parents bought me my own pc when I was 13, and ever since then, my love for videogames was cemented."

# Colour the single space between "mean" and "very very serious" red,
# matching the author's highlighted edit.
$spaceRange = $d.Content
if ($spaceRange.Find.Execute("I of course mean very")) {
    $meanEnd = $spaceRange.Start + "I of course mean".Length
    $spaceOnly = $d.Range($meanEnd, $meanEnd + 1)
    Write-Host "space char:[$($spaceOnly.Text)]"
    $spaceOnly.Font.Color = 255
}

# --- Paragraph 3 (handball) ---
Replace-Text "If I wasn't spending time at home playing videogames with my father, I was spending in the local sports hall playing handball, I eventually" "If I wasn't at home playing videogames with my father, I was spending time in the local sports hall playing handball. I eventually"
Replace-Text "with all my friends around, we'd have a lot of fun. I didn't enjoy the training as much as I though I would," "with all my friends around, we would have a lot of fun, but I also learned how important morale and teamwork is. I didn't enjoy the training as much as I thought I would,"
Replace-Text "While half of my days was spent playing handball" "While half of my days were spent playing handball"

# --- Paragraph 4 (flower industry) ---
Replace-Text "my love for programming started, this was only a seed though," "my love for programming started. This was only a seed though,"
Replace-Text "in the flower industry(Pun intended). While doing this job, I was torn" "in the flower industry (Pun intended). While doing this job I learned some essential skills like leadership skills & intercultural communication and I used those skills every day. But I was torn"
